# Update header labels on every worksheet:
#   A1: "Input Sheet" -> "Car Name"
#   B1: "Value"       -> "Values"

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").Value = "Car Name"
    $ws.Range("B1").Value = "Values"
}
